$wb = $excel.ActiveWorkbook
$wsRush = $wb.Worksheets.Item("Rushing")
$wsRecv = $wb.Worksheets.Item("Receiving")

# ================= Rushing sheet: Week 15 totals =================
$wsRush.Range("C2").Value = 7
$wsRush.Range("E2").Value = 9

$wsRush.Range("E3").Value = 2

$wsRush.Range("C4").Value = 90
$wsRush.Range("D4").Value = 59
$wsRush.Range("E4").Value = 14

$wsRush.Range("C8").Value = 12
$wsRush.Range("D8").Value = 11
$wsRush.Range("E8").Value = 3
$wsRush.Range("F8").Value = 6

$wsRush.Range("C12").Value = 5

# ================= Receiving sheet: Week 15 totals =================

# Row 6: D.Johnson logged for the first time this season
$wsRecv.Range("B6").Value = "D.Johnson"
$wsRecv.Range("C6").Value = 0
$wsRecv.Range("D6").Value = 0
$wsRecv.Range("E6").Value = 1
$wsRecv.Range("F6").Value = 1
$wsRecv.Range("G6").Value = 0
$wsRecv.Range("H6").Value = 0

# Row 7: D.Parker's updated season totals through Week 15
$wsRecv.Range("B7").Value = "D.Parker"
$wsRecv.Range("C7").Value = 77
$wsRecv.Range("D7").Value = 67
$wsRecv.Range("E7").Value = 15
$wsRecv.Range("F7").Value = 8
$wsRecv.Range("G7").Value = 7
$wsRecv.Range("H7").Value = 5

# Row 10: A.Wilson
$wsRecv.Range("C10").Value = 36
$wsRecv.Range("D10").Value = 24
$wsRecv.Range("E10").Value = 3

# Row 11: M.Hollins
$wsRecv.Range("C11").Value = 17
$wsRecv.Range("G11").Value = 7

# Row 12: I.Ford
$wsRecv.Range("C12").Value = 10
$wsRecv.Range("D12").Value = 8
$wsRecv.Range("E12").Value = 2
$wsRecv.Range("F12").Value = 2

# Row 14: M.Gesicki
$wsRecv.Range("C14").Value = 78
$wsRecv.Range("D14").Value = 52
$wsRecv.Range("E14").Value = 20
$wsRecv.Range("F14").Value = 12
$wsRecv.Range("G14").Value = 9
$wsRecv.Range("H14").Value = 7

# Row 16: A.Shaheen
$wsRecv.Range("C16").Value = 3

# Row 17: H.Long
$wsRecv.Range("C17").Value = 27
$wsRecv.Range("D17").Value = 22
$wsRecv.Range("G17").Value = 7
$wsRecv.Range("H17").Value = 4

# ================= Receiving sheet: Week 16 simulated =================
# J.Waddle explodes for a big game and is logged as a new row at the bottom
$wsRecv.Range("A18").Value = 16
$wsRecv.Range("B18").Value = "J.Waddle"
$wsRecv.Range("C18").Value = 108
$wsRecv.Range("D18").Value = 88
$wsRecv.Range("E18").Value = 19
$wsRecv.Range("F18").Value = 9
$wsRecv.Range("G18").Value = 15
$wsRecv.Range("H18").Value = 12

# Formatting for the new row's index cell: bold, centered, left/right border
$wsRecv.Range("A18").Font.Bold = $true
$wsRecv.Range("A18").HorizontalAlignment = -4108
$wsRecv.Range("A18").VerticalAlignment = -4160
$wsRecv.Range("A18").Borders.Item(7).LineStyle = 1
$wsRecv.Range("A18").Borders.Item(10).LineStyle = 1

$wsRecv.Range("H19").Select()
